# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The error-table rows represent forecast horizons (Q0..Q9, tracked via the
# "N" count in column G). A new leading observation is inserted at the top
# of the table (row 2): every existing row's ME/MAE/MSE/RMSE/SE values
# (columns B:F) shift down by one row, each row's sample count (column G)
# increases by one, and row 2 receives the newly computed error metrics for
# the additional observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original B:F values (rows 2-10) before anything is overwritten,
# since each row is about to receive the values currently held by the row
# above it.
$original = @{}
for ($r = 2; $r -le 10; $r++) {
    $original[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2
    )
}

# Every existing row's sample count (column G) grows by one observation.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 7).Value2 + 1
}

# Rows 3-11 inherit the B:F metrics that used to belong to the row above
# them (process bottom-up so sources aren't clobbered before being read).
for ($r = 11; $r -ge 3; $r--) {
    $prevValues = $original[$r - 1]
    $ws.Cells.Item($r, 2).Value = $prevValues[0]
    $ws.Cells.Item($r, 3).Value = $prevValues[1]
    $ws.Cells.Item($r, 4).Value = $prevValues[2]
    $ws.Cells.Item($r, 5).Value = $prevValues[3]
    $ws.Cells.Item($r, 6).Value = $prevValues[4]
}

# Row 2 (the new leading observation) gets the freshly computed error
# metrics from the corrected naive component forecaster.
$ws.Cells.Item(2, 2).Value = 0.1660323612164193
$ws.Cells.Item(2, 3).Value = 0.468872571632237
$ws.Cells.Item(2, 4).Value = 0.3770675980855352
$ws.Cells.Item(2, 5).Value = 0.6140583018619121
$ws.Cells.Item(2, 6).Value = 0.6119356173368786
